$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the team-specific transition/time-distribution matrix values
# (Houston_A) with newly computed figures. Each assignment below sets a
# single cell in the matrix to its updated value, per the source data
# update.
$ws.Cells.Item(2, 2).Value = 0.2163265306122449
$ws.Cells.Item(2, 3).Value = 0.5163265306122449
$ws.Cells.Item(2, 10).Value = 0.0163265306122449
$ws.Cells.Item(2, 15).Value = 0.001020408163265306
$ws.Cells.Item(2, 16).Value = 0.1571428571428571
$ws.Cells.Item(2, 19).Value = 0.09285714285714286
$ws.Cells.Item(3, 2).Value = 0.007561436672967864
$ws.Cells.Item(3, 3).Value = 0.02079395085066163
$ws.Cells.Item(3, 10).Value = 0.03213610586011342
$ws.Cells.Item(3, 16).Value = 0.7126654064272212
$ws.Cells.Item(3, 19).Value = 0.2268431001890359
$ws.Cells.Item(4, 10).Value = 0.06766917293233082
$ws.Cells.Item(4, 15).Value = 0.007518796992481203
$ws.Cells.Item(4, 16).Value = 0.5864661654135338
$ws.Cells.Item(4, 19).Value = 0.3383458646616541
$ws.Cells.Item(6, 2).Value = 0.06583072100313479
$ws.Cells.Item(6, 4).Value = 0.0109717868338558
$ws.Cells.Item(6, 6).Value = 0.054858934169279
$ws.Cells.Item(6, 10).Value = 0.2884012539184953
$ws.Cells.Item(6, 15).Value = 0.01567398119122257
$ws.Cells.Item(6, 17).Value = 0.1457680250783699
$ws.Cells.Item(6, 18).Value = 0.04858934169278997
$ws.Cells.Item(6, 19).Value = 0.3699059561128527
$ws.Cells.Item(7, 2).Value = 0.1066176470588235
$ws.Cells.Item(7, 4).Value = 0.02573529411764706
$ws.Cells.Item(7, 6).Value = 0.04227941176470588
$ws.Cells.Item(7, 10).Value = 0.1617647058823529
$ws.Cells.Item(7, 15).Value = 0.01470588235294118
$ws.Cells.Item(7, 17).Value = 0.1819852941176471
$ws.Cells.Item(7, 18).Value = 0.08272058823529412
$ws.Cells.Item(7, 19).Value = 0.3841911764705883
$ws.Cells.Item(8, 2).Value = 0.1031626506024096
$ws.Cells.Item(8, 4).Value = 0.01506024096385542
$ws.Cells.Item(8, 6).Value = 0.0572289156626506
$ws.Cells.Item(8, 10).Value = 0.1385542168674699
$ws.Cells.Item(8, 15).Value = 0.01355421686746988
$ws.Cells.Item(8, 17).Value = 0.1626506024096386
$ws.Cells.Item(8, 18).Value = 0.08207831325301204
$ws.Cells.Item(8, 19).Value = 0.427710843373494
$ws.Cells.Item(9, 2).Value = 0.1194852941176471
$ws.Cells.Item(9, 4).Value = 0.01102941176470588
$ws.Cells.Item(9, 6).Value = 0.05330882352941176
$ws.Cells.Item(9, 10).Value = 0.1470588235294118
$ws.Cells.Item(9, 15).Value = 0.02205882352941177
$ws.Cells.Item(9, 17).Value = 0.1452205882352941
$ws.Cells.Item(9, 18).Value = 0.08455882352941177
$ws.Cells.Item(9, 19).Value = 0.4172794117647059
$ws.Cells.Item(10, 2).Value = 0.1105710814094775
$ws.Cells.Item(10, 4).Value = 0.02162818955042527
$ws.Cells.Item(10, 5).Value = 0.0009720534629404617
$ws.Cells.Item(10, 6).Value = 0.06366950182260024
$ws.Cells.Item(10, 10).Value = 0.1360874848116646
$ws.Cells.Item(10, 15).Value = 0.01555285540704739
$ws.Cells.Item(10, 17).Value = 0.2240583232077764
$ws.Cells.Item(10, 18).Value = 0.07363304981773998
$ws.Cells.Item(10, 19).Value = 0.3538274605103281
$ws.Cells.Item(11, 7).Value = 0.1451612903225807
$ws.Cells.Item(11, 10).Value = 0.1075268817204301
$ws.Cells.Item(11, 11).Value = 0.2096774193548387
$ws.Cells.Item(11, 12).Value = 0.5204301075268817
$ws.Cells.Item(11, 19).Value = 0.01720430107526882
$ws.Cells.Item(12, 7).Value = 0.6921529175050302
$ws.Cells.Item(12, 10).Value = 0.2414486921529175
$ws.Cells.Item(12, 11).Value = 0.008048289738430584
$ws.Cells.Item(12, 12).Value = 0.03420523138832998
$ws.Cells.Item(12, 19).Value = 0.02414486921529175
$ws.Cells.Item(13, 7).Value = 0.6814159292035398
$ws.Cells.Item(13, 10).Value = 0.2212389380530974
$ws.Cells.Item(13, 19).Value = 0.09734513274336283
$ws.Cells.Item(14, 6).Value = 0.1111111111111111
$ws.Cells.Item(14, 7).Value = 0.5555555555555556
$ws.Cells.Item(14, 10).Value = 0.3333333333333333
$ws.Cells.Item(15, 6).Value = 0.02903225806451613
$ws.Cells.Item(15, 8).Value = 0.1290322580645161
$ws.Cells.Item(15, 9).Value = 0.07258064516129033
$ws.Cells.Item(15, 10).Value = 0.3854838709677419
$ws.Cells.Item(15, 11).Value = 0.07258064516129033
$ws.Cells.Item(15, 13).Value = 0.01451612903225807
$ws.Cells.Item(15, 14).Value = 0.001612903225806452
$ws.Cells.Item(15, 15).Value = 0.06290322580645161
$ws.Cells.Item(15, 19).Value = 0.232258064516129
$ws.Cells.Item(16, 6).Value = 0.02054794520547945
$ws.Cells.Item(16, 8).Value = 0.160958904109589
$ws.Cells.Item(16, 9).Value = 0.0684931506849315
$ws.Cells.Item(16, 10).Value = 0.4589041095890411
$ws.Cells.Item(16, 11).Value = 0.08904109589041095
$ws.Cells.Item(16, 13).Value = 0.0136986301369863
$ws.Cells.Item(16, 15).Value = 0.04794520547945205
$ws.Cells.Item(16, 19).Value = 0.1404109589041096
$ws.Cells.Item(17, 6).Value = 0.02144388849177984
$ws.Cells.Item(17, 8).Value = 0.1672623302358828
$ws.Cells.Item(17, 9).Value = 0.1007862759113653
$ws.Cells.Item(17, 10).Value = 0.4310221586847748
$ws.Cells.Item(17, 11).Value = 0.1029306647605432
$ws.Cells.Item(17, 13).Value = 0.01286633309506791
$ws.Cells.Item(17, 14).Value = 0.0007147962830593281
$ws.Cells.Item(17, 15).Value = 0.05289492494639028
$ws.Cells.Item(17, 19).Value = 0.1100786275911365
$ws.Cells.Item(18, 6).Value = 0.01879699248120301
$ws.Cells.Item(18, 8).Value = 0.1578947368421053
$ws.Cells.Item(18, 9).Value = 0.07142857142857142
$ws.Cells.Item(18, 10).Value = 0.4642857142857143
$ws.Cells.Item(18, 11).Value = 0.09210526315789473
$ws.Cells.Item(18, 13).Value = 0.01503759398496241
$ws.Cells.Item(18, 14).Value = 0.003759398496240601
$ws.Cells.Item(18, 15).Value = 0.06390977443609022
$ws.Cells.Item(18, 19).Value = 0.112781954887218
$ws.Cells.Item(19, 6).Value = 0.01428571428571429
$ws.Cells.Item(19, 8).Value = 0.2187012987012987
$ws.Cells.Item(19, 9).Value = 0.07350649350649351
$ws.Cells.Item(19, 10).Value = 0.3698701298701298
$ws.Cells.Item(19, 11).Value = 0.1124675324675325
$ws.Cells.Item(19, 13).Value = 0.02051948051948052
$ws.Cells.Item(19, 14).Value = 0.001818181818181818
$ws.Cells.Item(19, 15).Value = 0.06805194805194806
$ws.Cells.Item(19, 19).Value = 0.1207792207792208
